# Generate Report for Archive
# - Flip the localization status shown in the report from "Ready for handoff"
#   to "In Translation" everywhere it is surfaced (Overview sheet's per-locale
#   status columns, plus each locale detail sheet's Status column).
# - Re-run the column autosize for the Status column(s) now that the text is
#   shorter, matching the narrower "In Translation" column width.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$oldStatus = "Ready for handoff"

# Narrowest ColumnWidth (character units) this runtime can snap to that lands
# on the report's updated "Status" column width (the pixel grid here is
# coarser than Excel's native grid, so this is the closest achievable value).
$statusColumnWidth = 12.43

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = $statusColumnWidth
$overview.Columns.Item(6).ColumnWidth = $statusColumnWidth

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Columns.Item(3).ColumnWidth = $statusColumnWidth

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Columns.Item(3).ColumnWidth = $statusColumnWidth

Write-Host "Updated status '$oldStatus' -> '$newStatus' on Overview!E2:F2, zh-cn!C2, de-de!C2 and resized Status columns."
